$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1581.2778
$ws.Range("I19").Value = 1611.25
$ws.Range("J19").Value = 1341.5
$ws.Range("K19").Value = 1611.25
$ws.Range("L19").Value = 1341.5
$ws.Range("M19").Value = -1436.25
$ws.Range("N19").Value = -1691.5
$ws.Range("H33").Value = 376.73685
$ws.Range("I33").Value = 215.29411
$ws.Range("K33").Value = 215.29411
$ws.Range("M33").Value = 13.70589000000001
$ws.Range("H112").Value = 1963708
$ws.Range("I112").Value = 1900
$ws.Range("J112").Value = 2384095.2
$ws.Range("K112").Value = 5700
$ws.Range("L112").Value = 7152285.600000001
$ws.Range("M112").Value = -4592
$ws.Range("N112").Value = -7154501.600000001
$ws.Range("H132").Value = 1397.7462
$ws.Range("I132").Value = 1396.3966
$ws.Range("K132").Value = 4189.1898
$ws.Range("M132").Value = -1659.1898
$ws.Range("H138").Value = 2519.26
$ws.Range("I138").Value = 1270.6428
$ws.Range("K138").Value = 3811.9284
$ws.Range("M138").Value = 1328.0716

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31273968
$ws.Range("I32").Value = 36380452
$ws.Range("K32").Value = 36380452
$ws.Range("M32").Value = -36380165
$ws.Range("H43").Value = 24817.5
$ws.Range("J43").Value = 28033.334
$ws.Range("L43").Value = 28033.334
$ws.Range("N43").Value = -28659.334
$ws.Range("H61").Value = 3407.5386
$ws.Range("I61").Value = 3178.2856
$ws.Range("K61").Value = 3178.2856
$ws.Range("M61").Value = -2966.2856
$ws.Range("H112").Value = 85499.5
$ws.Range("J112").Value = 85499.5
$ws.Range("L112").Value = 85499.5
$ws.Range("N112").Value = -88453.5
$ws.Range("H114").Value = 59994.5
$ws.Range("J114").Value = 59994.5
$ws.Range("L114").Value = 59994.5
$ws.Range("N114").Value = -68672.5
$ws.Range("H122").Value = 3798.8
$ws.Range("I122").Value = 3582.3333
$ws.Range("J122").Value = 4123.5
$ws.Range("K122").Value = 10746.9999
$ws.Range("L122").Value = 12370.5
$ws.Range("M122").Value = -8296.999899999999
$ws.Range("N122").Value = -17270.5
$ws.Range("H136").Value = 3407.5386
$ws.Range("I136").Value = 3178.2856
$ws.Range("K136").Value = 9534.856800000001
$ws.Range("M136").Value = -6984.856800000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 8023
$ws.Range("I31").Value = 8023
$ws.Range("K31").Value = 8023
$ws.Range("M31").Value = -7771
$ws.Range("H99").Value = 3500
$ws.Range("J99").Value = 3800
$ws.Range("L99").Value = 3800
$ws.Range("N99").Value = -6796

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 260.35715
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 5028.1333
$ws.Range("I31").Value = 2585.3
$ws.Range("J31").Value = 6249.55
$ws.Range("K31").Value = 2585.3
$ws.Range("L31").Value = 6249.55
$ws.Range("M31").Value = -2290.3
$ws.Range("N31").Value = -6839.55
$ws.Range("H34").Value = 5028.1333
$ws.Range("I34").Value = 2585.3
$ws.Range("J34").Value = 6249.55
$ws.Range("K34").Value = 2585.3
$ws.Range("L34").Value = 6249.55
$ws.Range("M34").Value = -2383.3
$ws.Range("N34").Value = -6653.55
$ws.Range("H100").Value = 61701.75
$ws.Range("J100").Value = 61701.75
$ws.Range("L100").Value = 61701.75
$ws.Range("N100").Value = -63865.75
$ws.Range("H114").Value = 55839.5
$ws.Range("J114").Value = 55839.5
$ws.Range("L114").Value = 55839.5
$ws.Range("N114").Value = -64517.5
$ws.Range("H116").Value = 108653.5
$ws.Range("J116").Value = 108653.5
$ws.Range("L116").Value = 108653.5
$ws.Range("N116").Value = -117831.5
$ws.Range("H129").Value = 78046.8
$ws.Range("J129").Value = 90058.5
$ws.Range("L129").Value = 90058.5
$ws.Range("N129").Value = -100058.5
$ws.Range("H134").Value = 3061.158
$ws.Range("I134").Value = 3015.75
$ws.Range("J134").Value = 3303.3333
$ws.Range("K134").Value = 9047.25
$ws.Range("L134").Value = 9909.999899999999
$ws.Range("M134").Value = -6512.25
$ws.Range("N134").Value = -14979.9999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6106
$ws.Range("I3").Value = 6106
$ws.Range("K3").Value = 18318
$ws.Range("M3").Value = -18206
$ws.Range("H95").Value = 7500
$ws.Range("J95").Value = 7500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -26618
$ws.Range("H131").Value = 1779.875
$ws.Range("J131").Value = 1754.1333
$ws.Range("L131").Value = 5262.3999
$ws.Range("N131").Value = -15342.3999
$ws.Range("H136").Value = 2752.4167
$ws.Range("I136").Value = 2128.625
$ws.Range("K136").Value = 6385.875
$ws.Range("M136").Value = -1285.875
$ws.Range("H137").Value = 5513.7144
$ws.Range("I137").Value = 1999.5
$ws.Range("J137").Value = 6919.4
$ws.Range("K137").Value = 5998.5
$ws.Range("L137").Value = 20758.2
$ws.Range("M137").Value = -898.5
$ws.Range("N137").Value = -30958.2

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 50045
$ws.Range("I2").Value = 90
$ws.Range("K2").Value = 90
$ws.Range("M2").Value = 23

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3700.3125
$ws.Range("I40").Value = 2823.4614
$ws.Range("K40").Value = 2823.4614
$ws.Range("M40").Value = -2687.4614
$ws.Range("H46").Value = 3456.2778
$ws.Range("I46").Value = 449.33334
$ws.Range("J46").Value = 4057.6667
$ws.Range("K46").Value = 449.33334
$ws.Range("L46").Value = 4057.6667
$ws.Range("M46").Value = -261.33334
$ws.Range("N46").Value = -4433.6667
$ws.Range("H110").Value = 40001
$ws.Range("J110").Value = 40001
$ws.Range("L110").Value = 40001
$ws.Range("N110").Value = -48181
$ws.Range("H132").Value = 479276.8
$ws.Range("I132").Value = 717237.9
$ws.Range("J132").Value = 3354.7144
$ws.Range("K132").Value = 2151713.7
$ws.Range("L132").Value = 10064.1432
$ws.Range("M132").Value = -2149183.7
$ws.Range("N132").Value = -15124.1432

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 86401
$ws.Range("J103").Value = 86401
$ws.Range("L103").Value = 86401
$ws.Range("N103").Value = -88745
$ws.Range("H132").Value = 52855.75
$ws.Range("I132").Value = 55469.26
$ws.Range("K132").Value = 166407.78
$ws.Range("M132").Value = -163877.78
